$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price entry was recorded. Insert a new row at row 59 (pushing the
# existing rows 59-63 down to 60-64) and populate it with the new observation.
$ws.Rows.Item(59).Insert()

$ws.Cells.Item(59, 1).Value = 11
$ws.Cells.Item(59, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(59, 3).Value = "Bíobío"
$ws.Cells.Item(59, 4).Value = 45021
$ws.Cells.Item(59, 5).Value = 8
$ws.Cells.Item(59, 6).Value = 100112043
$ws.Cells.Item(59, 7).Value = "Pepino dulce"
$ws.Cells.Item(59, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(59, 9).Value = "Primera"
$ws.Cells.Item(59, 10).Value = 100
$ws.Cells.Item(59, 11).Value = 13000
$ws.Cells.Item(59, 12).Value = 14000
$ws.Cells.Item(59, 13).Value = 13500
$ws.Cells.Item(59, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(59, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(59, 16).Value = 750
$ws.Cells.Item(59, 17).Value = 18
$ws.Cells.Item(59, 18).Value = "Hortaliza"
